$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the existing row 963, shifting old rows 963-1032 down to 967-1036
$ws.Rows("963:966").Insert()

$row963 = New-Object 'object[,]' 1,18
$row963[0,0] = 9
$row963[0,1] = 'Vega Central Mapocho de Santiago'
$row963[0,2] = 'Metropolitana'
$row963[0,3] = 45106
$row963[0,4] = 13
$row963[0,5] = 100112023
$row963[0,6] = 'Brócoli'
$row963[0,7] = 'Sin especificar'
$row963[0,8] = 'Primera'
$row963[0,9] = 3400
$row963[0,10] = 700
$row963[0,11] = 800
$row963[0,12] = 750
$row963[0,13] = '$/unidad'
$row963[0,14] = 'Región Metropolitana'
$row963[0,15] = 750
$row963[0,16] = 1
$row963[0,17] = 'Hortaliza'
$ws.Range("A963:R963").Value = $row963

$row964 = New-Object 'object[,]' 1,18
$row964[0,0] = 9
$row964[0,1] = 'Vega Central Mapocho de Santiago'
$row964[0,2] = 'Metropolitana'
$row964[0,3] = 45106
$row964[0,4] = 13
$row964[0,5] = 100112023
$row964[0,6] = 'Brócoli'
$row964[0,7] = 'Sin especificar'
$row964[0,8] = 'Primera'
$row964[0,9] = 1600
$row964[0,10] = 800
$row964[0,11] = 900
$row964[0,12] = 850
$row964[0,13] = '$/unidad'
$row964[0,14] = 'Región de O''Higgins'
$row964[0,15] = 850
$row964[0,16] = 1
$row964[0,17] = 'Hortaliza'
$ws.Range("A964:R964").Value = $row964

$row965 = New-Object 'object[,]' 1,18
$row965[0,0] = 9
$row965[0,1] = 'Vega Central Mapocho de Santiago'
$row965[0,2] = 'Metropolitana'
$row965[0,3] = 45106
$row965[0,4] = 13
$row965[0,5] = 100112023
$row965[0,6] = 'Brócoli'
$row965[0,7] = 'Sin especificar'
$row965[0,8] = 'Segunda'
$row965[0,9] = 1690
$row965[0,10] = 600
$row965[0,11] = 600
$row965[0,12] = 600
$row965[0,13] = '$/unidad'
$row965[0,14] = 'Región Metropolitana'
$row965[0,15] = 600
$row965[0,16] = 1
$row965[0,17] = 'Hortaliza'
$ws.Range("A965:R965").Value = $row965

$row966 = New-Object 'object[,]' 1,18
$row966[0,0] = 9
$row966[0,1] = 'Vega Central Mapocho de Santiago'
$row966[0,2] = 'Metropolitana'
$row966[0,3] = 45106
$row966[0,4] = 13
$row966[0,5] = 100112023
$row966[0,6] = 'Brócoli'
$row966[0,7] = 'Sin especificar'
$row966[0,8] = 'Segunda'
$row966[0,9] = 970
$row966[0,10] = 700
$row966[0,11] = 700
$row966[0,12] = 700
$row966[0,13] = '$/unidad'
$row966[0,14] = 'Región de O''Higgins'
$row966[0,15] = 700
$row966[0,16] = 1
$row966[0,17] = 'Hortaliza'
$ws.Range("A966:R966").Value = $row966
